$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the new text value for cell E12 (row 12 in "Сельское хозяйство" block)
$ws.Range("E12").Value = "Урожайность овощей - harvest (цент.) (8007025)"

# Update the selection to reflect the cell that was last active (E17)
$ws.Range("E17").Select()
